$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "白细胞"
$ws.Range("C2").Value = "4-10"

$ws.Range("A3").Value = "中性粒细胞百分率"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "46.7"
$ws.Range("C3").Value = "50-70"

$ws.Range("A4").Value = "淋巴细胞百分率"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "41.6"
$ws.Range("C4").Value = "20-40"

$ws.Range("A5").Value = "单核细胞百分率"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "5.7"
$ws.Range("C5").Value = "3-8"

$ws.Range("A6").Value = "嗜酸性粒细胞百分率"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "5.6"
$ws.Range("C6").Value = "0.5-5"

$ws.Range("A7").Value = "嗜碱性粒细胞百分率"
$ws.Range("C7").Value = "0-1"

$ws.Range("A8").Value = "中性粒细胞绝对值"
$ws.Range("C8").Value = "2-7"

$ws.Range("A9").Value = "淋巴细胞绝对值"
$ws.Range("C9").Value = "0.8-4"

$ws.Range("A10").Value = "单核细胞绝对值"
$ws.Range("C10").Value = "0.1-1"

$ws.Range("A11").Value = "嗜酸性粒细胞绝对值"
$ws.Range("C11").Value = "0.05-0.5"

$ws.Range("A12").Value = "嗜碱性粒细胞绝对值"
$ws.Range("C12").Value = "0-0.1"

$ws.Range("A13").Value = "红细胞"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "5.2"
$ws.Range("C13").Value = "4-5.5"

$ws.Range("A14").Value = "血红蛋白"
$ws.Range("C14").Value = "110-160"

$ws.Range("A15").Value = "红细胞压积"
$ws.Range("C15").Value = "32-50"

$ws.Range("A16").Value = "平均红细胞体积"
$ws.Range("C16").Value = "80-97"

$ws.Range("A17").Value = "RBC平均HB含量"
$ws.Range("C17").Value = "26-32"

$ws.Range("A18").Value = "RBC平均HB浓度"
$ws.Range("C18").Value = "320-360"

$ws.Range("A19").Value = "红细胞分布宽度"
$ws.Range("C19").Value = "10-15"

$ws.Range("A20").Value = "红细胞分布宽度"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "50.6"
$ws.Range("C20").Value = "39-46"

$ws.Range("A21").Value = "血小板"
$ws.Range("C21").Value = "100-300"

$ws.Range("A22").Value = "平均血小板体积"
$ws.Range("C22").Value = "9.4-12.5"

$ws.Range("A23").Value = "血小板分布宽度"
$ws.Range("C23").Value = "15.5-18.1"

$ws.Range("A24").Value = "血小板压积"
$ws.Range("C24").Value = "0.108-0.282"
